# "aggiornamento fino a 8/12" - append daily data rows from 2021-09-21
# (serial 44460) through 2021-12-08 (serial 44538), i.e. worksheet rows
# 386-464, extending the dimension from A1:D385 to A1:D464.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, replicate the formatting of the last existing data row (385)
# across the whole new block, so column A keeps its date number format,
# center/top alignment and thin border (style index 2), while columns
# B:D keep the default (unstyled) look - exactly like all prior rows.
$ws.Range("A385:D385").Copy()
$ws.Range("A386:D464").PasteSpecial(-4122)

# Rows 386-460: plain continuation, all zeros, one calendar day per row.
$serial = 44460
for ($r = 386; $r -le 460; $r++) {
    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $serial = $serial + 1
}

# Rows 461-464: new positive cases starting 2021-12-05 (serial 44535).
$ws.Cells.Item(461, 1).Value = 44535
$ws.Cells.Item(461, 2).Value = 1
$ws.Cells.Item(461, 3).Value = 1
$ws.Cells.Item(461, 4).Value = 145.7725947521866

$ws.Cells.Item(462, 1).Value = 44536
$ws.Cells.Item(462, 2).Value = 0
$ws.Cells.Item(462, 3).Value = 1
$ws.Cells.Item(462, 4).Value = 145.7725947521866

$ws.Cells.Item(463, 1).Value = 44537
$ws.Cells.Item(463, 2).Value = 0
$ws.Cells.Item(463, 3).Value = 1
$ws.Cells.Item(463, 4).Value = 145.7725947521866

$ws.Cells.Item(464, 1).Value = 44538
$ws.Cells.Item(464, 2).Value = 0
$ws.Cells.Item(464, 3).Value = 1
$ws.Cells.Item(464, 4).Value = 145.7725947521866
